$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "2025-04-26 21:04:19"
$ws.Range("C25").Value = "John Smith moved battery 7 from No Location to floor space 1.`nNow John Smith is Confident.`n"
$ws.Range("C25").WrapText = $true

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "2025-04-26 21:09:16"
$ws.Range("C26").Value = "John Smith took picture of battery 7.`nNow John Smith is Happy.`n"
$ws.Range("C26").WrapText = $true

$ws.Rows.Item(25).EntireRow.AutoFit()
$ws.Rows.Item(26).EntireRow.AutoFit()
